$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.814.62'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '1.607.93'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'210.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.54%  '
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'0.249"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').Value = "'19.74"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.73%  '
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D12').Value = '1.832.72'
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '1.601.30'
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('D16').Value = '26.780.21'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').Value = "'63.70"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.00%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = "'210.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = "'6.75"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('E22').Value = '  -2.25%  '
$ws.Range('D23').Value = "'2.33"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.74%  '
$ws.Range('D24').Value = "'8.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.46%  '
$ws.Range('D25').Value = "'146.53"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('E28').Value = '  -4.22%  '
$ws.Range('E29').Value = '  -1.43%  '
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E31').Value = '  -2.51%  '
$ws.Range('E32').Value = '  -2.75%  '
$ws.Range('E33').Value = '  +24.61%  '
$ws.Range('D35').Value = '1.318.86'
$ws.Range('E35').Value = '  +1.42%  '
$ws.Range('E36').Value = '  -2.00%  '
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('D39').Value = "'0.822"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').Value = "'2.20"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.34%  '
$ws.Range('D43').Value = "'5.28"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').Value = "'62.97"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.30%  '
$ws.Range('D45').Value = '1.744.89'
$ws.Range('E45').Value = '  -1.17%  '
$ws.Range('D46').Value = "'89.08"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('E47').Value = '  +1.20%  '
$ws.Range('D48').Value = "'0.818"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.81%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  -4.37%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.0510"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.0980"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.13%  '
